# Updates the cryptocurrency price/volume table (GitHub Actions scheduled refresh).
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row layout: Row, Coin, Link, Price, "Volume(1h)"
$data = @(
    ,@(2, 'Bitcoin', 'https://coinranking.com/coin/Qwsogvtv82FCd+bitcoin-btc', '29.376.80', '  -0.02%  ')
    ,@(3, 'Ethereum', 'https://coinranking.com/coin/razxDUgYGNAdQ+ethereum-eth', '1.846.14', '  -0.14%  ')
    ,@(4, 'TetherUSD', 'https://coinranking.com/coin/HIVsRcGKkPFtW+tetherusd-usdt', '0.9979', '  -0.13%  ')
    ,@(5, 'BNB', 'https://coinranking.com/coin/WcwrkfNI4FUAe+bnb-bnb', '240.50', '  +0.10%  ')
    ,@(6, 'XRP', 'https://coinranking.com/coin/-l8Mn2pVlRs-p+xrp-xrp', '0.6272', '  -0.09%  ')
    ,@(7, 'USDC', 'https://coinranking.com/coin/aKzUVe4Hh_CON+usdc-usdc', '0.9998', '  +0.00%  ')
    ,@(8, 'Dogecoin', 'https://coinranking.com/coin/a91GCGd_u96cF+dogecoin-doge', '0.07497', '  -1.56%  ')
    ,@(9, 'Cardano', 'https://coinranking.com/coin/qzawljRxB5bYu+cardano-ada', '0.2901', '  -0.16%  ')
    ,@(10, 'Solana', 'https://coinranking.com/coin/zNZHO_Sjf+solana-sol', '24.49', '  -1.04%  ')
    ,@(11, 'TRON', 'https://coinranking.com/coin/qUhEFk1I61atv+tron-trx', '0.07733', '  -0.07%  ')
    ,@(12, 'WrappedEther', 'https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth', '1.846.26', '  -2.22%  ')
    ,@(13, 'Polkadot', 'https://coinranking.com/coin/25W7FG7om+polkadot-dot', '4.998', '  -0.66%  ')
    ,@(14, 'Polygon', 'https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic', '0.6815', '  +0.30%  ')
    ,@(15, 'ShibaInu', 'https://coinranking.com/coin/xz24e0BjL+shibainu-shib', '0.00001057', '  +0.13%  ')
    ,@(16, 'Litecoin', 'https://coinranking.com/coin/D7B1x_ks7WhV5+litecoin-ltc', '82.12', '  -1.18%  ')
    ,@(17, 'Uniswap', 'https://coinranking.com/coin/_H5FVG9iW+uniswap-uni', '6.186', '  +0.41%  ')
    ,@(18, 'WrappedBTC', 'https://coinranking.com/coin/x4WXHge-vvFY+wrappedbtc-wbtc', '29.416.36', '  +0.03%  ')
    ,@(19, 'BitcoinCash', 'https://coinranking.com/coin/ZlZpzOJo43mIo+bitcoincash-bch', '229.46', '  +0.78%  ')
    ,@(20, 'Avalanche', 'https://coinranking.com/coin/dvUj0CzDZ+avalanche-avax', '12.32', '  -0.26%  ')
    ,@(21, 'Dai', 'https://coinranking.com/coin/MoTuySvg7+dai-dai', '0.9994', '  -0.02%  ')
    ,@(22, 'Chainlink', 'https://coinranking.com/coin/VLqpJwogdhHNb+chainlink-link', '7.487', '  +0.05%  ')
    ,@(23, 'BinanceUSD', 'https://coinranking.com/coin/vSo2fu9iE1s0Y+binanceusd-busd', '0.9994', '  +0.03%  ')
    ,@(24, 'Monero', 'https://coinranking.com/coin/3mVx2FX_iJFp5+monero-xmr', '159.41', '  +0.40%  ')
    ,@(25, 'Stellar', 'https://coinranking.com/coin/f3iaFeCKEmkaZ+stellar-xlm', '0.1373', '  -1.00%  ')
    ,@(26, 'Cosmos', 'https://coinranking.com/coin/Knsels4_Ol-Ny+cosmos-atom', '8.423', '  +0.18%  ')
    ,@(27, 'EthereumClassic', 'https://coinranking.com/coin/hnfQfsYfeIGUQ+ethereumclassic-etc', '17.53', '  -0.90%  ')
    ,@(28, 'Hedera', 'https://coinranking.com/coin/jad286TjB+hedera-hbar', '0.06499', '  +15.92%  ')
    ,@(29, 'Toncoin', 'https://coinranking.com/coin/67YlI0K1b+toncoin-ton', '1.419', '  +0.60%  ')
    ,@(30, 'PancakeSwap', 'https://coinranking.com/coin/ncYFcP709+pancakeswap-cake', '1.484', '  +1.55%  ')
    ,@(31, 'Filecoin', 'https://coinranking.com/coin/ymQub4fuB+filecoin-fil', '4.097', '  -0.38%  ')
    ,@(32, 'InternetComputer(DFINITY)', 'https://coinranking.com/coin/aMNLwaUbY+internetcomputerdfinity-icp', '4.097', '  +0.72%  ')
    ,@(33, 'LidoDAOToken', 'https://coinranking.com/coin/Pe93bIOD2+lidodaotoken-ldo', '1.833', '  -0.14%  ')
    ,@(34, 'ARBITRUM', 'https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb', '1.142', '  -1.86%  ')
    ,@(35, 'ImmutableX', 'https://coinranking.com/coin/Z96jIvLU7+immutablex-imx', '0.6986', '  +0.15%  ')
    ,@(36, 'HuobiToken', 'https://coinranking.com/coin/DXwP4wF9ksbBO+huobitoken-ht', '2.580', '  -0.33%  ')
    ,@(37, 'Maker', 'https://coinranking.com/coin/qFakph2rpuMOL+maker-mkr', '1.266.83', '  +2.96%  ')
    ,@(38, 'MXToken', 'https://coinranking.com/coin/QUC5kVAxSoB-+mxtoken-mx', '2.836', '  +3.95%  ')
    ,@(39, 'VeChain', 'https://coinranking.com/coin/FEbS54wxo4oIl+vechain-vet', '0.01830', '  +1.50%  ')
    ,@(40, 'FraxShare', 'https://coinranking.com/coin/3nNpuxHJ8+fraxshare-fxs', '6.763', '  +5.87%  ')
    ,@(41, 'TrustWalletToken', 'https://coinranking.com/coin/Hm3OlynlC+trustwallettoken-twt', '0.9100', '  +1.14%  ')
    ,@(42, 'PaxDollar', 'https://coinranking.com/coin/JCKLgWPAF+paxdollar-usdp', '0.9994', '  -0.06%  ')
    ,@(43, 'RocketPoolETH', 'https://coinranking.com/coin/QJZRUGyNI+rocketpooleth-reth', '2.007.47', '  -18.42%  ')
    ,@(44, 'Quant', 'https://coinranking.com/coin/bauj_21eYVwso+quant-qnt', '101.28', '  -0.20%  ')
    ,@(45, 'Aave', 'https://coinranking.com/coin/ixgUfzmLR+aave-aave', '66.32', '  +0.57%  ')
    ,@(46, 'RenderToken', 'https://coinranking.com/coin/7C4Mh4xy1yDel+rendertoken-rndr', '1.741', '  +3.89%  ')
    ,@(47, 'BabyDogeCoin', 'https://coinranking.com/coin/JY1_q2c0g+babydogecoin-babydoge', '0.00000000119', '  +0.11%  ')
    ,@(48, 'Aptos', 'https://coinranking.com/coin/HGYj5JCv5+aptos-apt', '7.078', '  -2.00%  ')
    ,@(49, 'Algorand', 'https://coinranking.com/coin/TpHE2IShQw-sJ+algorand-algo', '0.1165', '  +2.35%  ')
    ,@(50, 'EnergySwap', 'https://coinranking.com/coin/SbWqqTui-+energyswap-ens', '9.049', '  +0.49%  ')
    ,@(51, 'TheSandbox', 'https://coinranking.com/coin/pxtKbG5rg+thesandbox-sand', '0.3973', '  -0.73%  ')
)

foreach ($item in $data) {
    $row   = $item[0]
    $coin  = $item[1]
    $link  = $item[2]
    $price = $item[3]
    $vol   = $item[4]

    $ws.Cells.Item($row, 2).Value = $coin
    $ws.Cells.Item($row, 3).Value = $link

    # Price column holds numeric-looking strings (e.g. "0.9100", "240.50")
    # that must stay as text so trailing zeros / exact formatting survive -
    # force text format, assign, then drop back to the default cell style
    # so no stray formatting is left behind.
    $priceCell = $ws.Cells.Item($row, 4)
    $priceCell.NumberFormat = "@"
    $priceCell.Value = $price
    $priceCell.Style = "Normal"

    $ws.Cells.Item($row, 5).Value = $vol
}
